$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (previously Melbourne/Left Bank) -> McKinnon/Hotlocks
$ws.Range("A8").Value = "McKinnon"
$ws.Range("B8").Value = "Hotlocks By Rachael Hairdresser, 260 McKinnon Road, McKinnon VIC 3204"
$ws.Range("C8").Value = "23/12/20 4:00pm-6:00pm"
$ws.Range("D8").Value = "Case had hair cut in store"

# Update row 9 (previously Melbourne/Lion Hotel) -> Melbourne/Left Bank
$ws.Range("B9").Value = "Left Bank Melbourne, 1 Southbank Blvd"
$ws.Range("C9").Value = "25/12/20 12pm - 3pm"
$ws.Range("D9").Value = "Case ate in store"

# Add new row 12 -> Southbank/Rockpool Bar and Grill
$ws.Range("A12").Value = "Southbank"
$ws.Range("B12").Value = "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank"
$ws.Range("C12").Value = "23/12/20 8:00pm-11:00pm"
$ws.Range("D12").Value = "Case attended restaurant"
